$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 26: ID 24 - InstinctiveEvasion (PassiveSkill)
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A26:E26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "InstinctiveEvasion"
$ws.Range("C26").Value = "PassiveSkill"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0

# New row 27: ID 25 - MarkingDodge (PassiveSkill)
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A27:E27").PasteSpecial(-4122) | Out-Null
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "MarkingDodge"
$ws.Range("C27").Value = "PassiveSkill"
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0

$excel.CutCopyMode = 0

$ws.Range("F28").Select() | Out-Null
